{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same text substitutions described by the OOXML diff:\n//  - Title/H1 text + its bold/italic repeats at the end of the doc\n//  - \"What we like\" bullet list items\n//  - \"What we don't like\" bullet list items\n//  - The italic summary paragraph near the end\n\nconst replacements = [\n  [\n    \"Play Mystic Spirits for free - Review of Red Rake Gaming Slot Game\",\n    \"Play Mystic Spirits Free - Exciting Gameplay and Stunning Graphics\",\n  ],\n  [\n    \"Impressive 50 paylines for more chances to win\",\n    \"Impressive 50 paylines\",\n  ],\n  [\n    \"Stunning 3D graphics with a Native American vibe\",\n    \"Stunning graphics with visually stunning effects\",\n  ],\n  [\n    \"Free Spins feature can result in up to 300 Free Spins\",\n    \"Exciting special symbols and features\",\n  ],\n  [\n    \"Potential win of 750x the initial bet\",\n    \"Potential win of 750 times the initial bet\",\n  ],\n  [\n    \"The game's volatility is average\",\n    \"Average volatility\",\n  ],\n  [\n    \"The RTP is lower than some other online slot games\",\n    \"Minimum bet of \u20ac0.20\",\n  ],\n  [\n    \"Discover the Native American-themed Mystic Spirits slot game by Red Rake Gaming. Play for free and read our expert review with pros and cons, RTP, and features.\",\n    \"Discover Mystic Spirits, an exciting slot game with stunning graphics and exciting features. Play for free!\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same text substitutions described by the OOXML diff using\n# Find/Replace over the whole document content (wdReplaceAll), which also\n# correctly updates the two occurrences of the title/heading text.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Mystic Spirits for free - Review of Red Rake Gaming Slot Game\", \"Play Mystic Spirits Free - Exciting Gameplay and Stunning Graphics\"),\n    @(\"Impressive 50 paylines for more chances to win\", \"Impressive 50 paylines\"),\n    @(\"Stunning 3D graphics with a Native American vibe\", \"Stunning graphics with visually stunning effects\"),\n    @(\"Free Spins feature can result in up to 300 Free Spins\", \"Exciting special symbols and features\"),\n    @(\"Potential win of 750x the initial bet\", \"Potential win of 750 times the initial bet\"),\n    @(\"The game's volatility is average\", \"Average volatility\"),\n    @(\"The RTP is lower than some other online slot games\", \"Minimum bet of \u20ac0.20\"),\n    @(\"Discover the Native American-themed Mystic Spirits slot game by Red Rake Gaming. Play for free and read our expert review with pros and cons, RTP, and features.\", \"Discover Mystic Spirits, an exciting slot game with stunning graphics and exciting features. Play for free!\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
